{"js": "// Replace the date line and every \"a\u00f7b=c, r\" answer cell with the new\n// values from the day's worksheet. Each \"before\" string is unique in the\n// document, so an exact (non-wildcard) search-and-replace is safe.\nconst replacements = [\n  [\"2025-06-08 Sunday\", \"2025-06-09 Monday\"],\n  [\"91\u00f79=10, 1\", \"96\u00f72=48, 0\"],\n  [\"68\u00f76=11, 2\", \"14\u00f72=7, 0\"],\n  [\"81\u00f75=16, 1\", \"36\u00f73=12, 0\"],\n  [\"82\u00f74=20, 2\", \"84\u00f74=21, 0\"],\n  [\"51\u00f79=5, 6\", \"31\u00f72=15, 1\"],\n  [\"53\u00f75=10, 3\", \"63\u00f72=31, 1\"],\n  [\"10\u00f75=2, 0\", \"86\u00f75=17, 1\"],\n  [\"92\u00f76=15, 2\", \"79\u00f76=13, 1\"],\n  [\"15\u00f72=7, 1\", \"31\u00f77=4, 3\"],\n  [\"92\u00f72=46, 0\", \"19\u00f77=2, 5\"],\n  [\"96\u00f75=19, 1\", \"39\u00f72=19, 1\"],\n  [\"22\u00f74=5, 2\", \"67\u00f73=22, 1\"],\n  [\"37\u00f74=9, 1\", \"40\u00f78=5, 0\"],\n  [\"12\u00f72=6, 0\", \"88\u00f76=14, 4\"],\n  [\"69\u00f73=23, 0\", \"48\u00f74=12, 0\"],\n  [\"17\u00f78=2, 1\", \"96\u00f73=32, 0\"],\n  [\"90\u00f76=15, 0\", \"73\u00f79=8, 1\"],\n  [\"33\u00f77=4, 5\", \"87\u00f76=14, 3\"],\n  [\"56\u00f75=11, 1\", \"37\u00f77=5, 2\"],\n  [\"23\u00f72=11, 1\", \"30\u00f76=5, 0\"],\n  [\"22\u00f72=11, 0\", \"59\u00f73=19, 2\"],\n  [\"47\u00f72=23, 1\", \"11\u00f75=2, 1\"],\n  [\"97\u00f72=48, 1\", \"81\u00f73=27, 0\"],\n  [\"26\u00f74=6, 2\", \"76\u00f74=19, 0\"],\n  [\"46\u00f75=9, 1\", \"34\u00f72=17, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00f7b=c, r\" answer cell with the new\n# values from the day's worksheet. Each \"before\" string is unique in the\n# document, so Find/Replace (wdReplaceAll) on the whole document is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-08 Sunday\", \"2025-06-09 Monday\"),\n    @(\"91\u00f79=10, 1\", \"96\u00f72=48, 0\"),\n    @(\"68\u00f76=11, 2\", \"14\u00f72=7, 0\"),\n    @(\"81\u00f75=16, 1\", \"36\u00f73=12, 0\"),\n    @(\"82\u00f74=20, 2\", \"84\u00f74=21, 0\"),\n    @(\"51\u00f79=5, 6\", \"31\u00f72=15, 1\"),\n    @(\"53\u00f75=10, 3\", \"63\u00f72=31, 1\"),\n    @(\"10\u00f75=2, 0\", \"86\u00f75=17, 1\"),\n    @(\"92\u00f76=15, 2\", \"79\u00f76=13, 1\"),\n    @(\"15\u00f72=7, 1\", \"31\u00f77=4, 3\"),\n    @(\"92\u00f72=46, 0\", \"19\u00f77=2, 5\"),\n    @(\"96\u00f75=19, 1\", \"39\u00f72=19, 1\"),\n    @(\"22\u00f74=5, 2\", \"67\u00f73=22, 1\"),\n    @(\"37\u00f74=9, 1\", \"40\u00f78=5, 0\"),\n    @(\"12\u00f72=6, 0\", \"88\u00f76=14, 4\"),\n    @(\"69\u00f73=23, 0\", \"48\u00f74=12, 0\"),\n    @(\"17\u00f78=2, 1\", \"96\u00f73=32, 0\"),\n    @(\"90\u00f76=15, 0\", \"73\u00f79=8, 1\"),\n    @(\"33\u00f77=4, 5\", \"87\u00f76=14, 3\"),\n    @(\"56\u00f75=11, 1\", \"37\u00f77=5, 2\"),\n    @(\"23\u00f72=11, 1\", \"30\u00f76=5, 0\"),\n    @(\"22\u00f72=11, 0\", \"59\u00f73=19, 2\"),\n    @(\"47\u00f72=23, 1\", \"11\u00f75=2, 1\"),\n    @(\"97\u00f72=48, 1\", \"81\u00f73=27, 0\"),\n    @(\"26\u00f74=6, 2\", \"76\u00f74=19, 0\"),\n    @(\"46\u00f75=9, 1\", \"34\u00f72=17, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
